$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.823.02'
$ws.Range('E2').Value = '  +3.06%  '
$ws.Range('D3').Value = '2.536.85'
$ws.Range('E3').Value = '  +5.56%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.31'
$ws.Range('E5').Value = '  +2.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.44'
$ws.Range('E6').Value = '  +5.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.589'
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('D9').Value = '2.536.63'
$ws.Range('E9').Value = '  +5.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.106'
$ws.Range('E10').Value = '  +2.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.77'
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('E13').Value = '  +3.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.87'
$ws.Range('E14').Value = '  +8.61%  '
$ws.Range('D15').Value = '2.993.29'
$ws.Range('E15').Value = '  +5.68%  '
$ws.Range('D16').Value = '63.628.61'
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('E17').Value = '  +3.97%  '
$ws.Range('D18').Value = '2.541.02'
$ws.Range('E18').Value = '  +5.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.45'
$ws.Range('E19').Value = '  +4.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '343.34'
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.35'
$ws.Range('E21').Value = '  +3.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.88'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  +0.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.88'
$ws.Range('E24').Value = '  +1.68%  '
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('E26').Value = '  +5.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.24'
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.44'
$ws.Range('E29').Value = '  +4.71%  '
$ws.Range('D30').Value = '0.0₃0823'
$ws.Range('E30').Value = '  +7.52%  '
$ws.Range('E31').Value = '  +4.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.81'
$ws.Range('E32').Value = '  +7.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '176.84'
$ws.Range('E33').Value = '  +3.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.54'
$ws.Range('E34').Value = '  +9.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '416.15'
$ws.Range('E35').Value = '  +16.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.402'
$ws.Range('E36').Value = '  +2.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.11'
$ws.Range('E37').Value = '  +3.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.42'
$ws.Range('E38').Value = '  -2.44%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.75'
$ws.Range('E40').Value = '  +5.42%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '40.93'
$ws.Range('E42').Value = '  +5.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '153.02'
$ws.Range('E43').Value = '  +6.53%  '
$ws.Range('E44').Value = '  +3.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.99'
$ws.Range('E45').Value = '  +2.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.610'
$ws.Range('E46').Value = '  +4.68%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0530'
$ws.Range('E47').Value = '  +2.55%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0967'
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.01'
$ws.Range('E49').Value = '  +6.98%  '
$ws.Range('E50').Value = '  +4.75%  '
$ws.Range('D51').Value = '0.0₆0231'
$ws.Range('E51').Value = '  +7.31%  '
